# Weekly update: insert a new "Zanahoria" price record for Vega Modelo de
# Temuco at row 168, pushing the existing rows 168-190 down to 169-191.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 168 (this shifts rows 168..190 to
# 169..191 and carries their formatting/styles with them, including the
# date number format on column D).
$ws.Rows.Item(168).Insert()

# Populate the newly inserted row 168 with the new weekly record.
$ws.Cells.Item(168, 1).Value = 10
$ws.Cells.Item(168, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(168, 3).Value = "La Araucanía"
$ws.Cells.Item(168, 4).Value = 44449
$ws.Cells.Item(168, 5).Value = 9
$ws.Cells.Item(168, 6).Value = 100114013
$ws.Cells.Item(168, 7).Value = "Zanahoria"
$ws.Cells.Item(168, 8).Value = "Sin especificar"
$ws.Cells.Item(168, 9).Value = "Primera"
$ws.Cells.Item(168, 10).Value = 125
$ws.Cells.Item(168, 11).Value = 5500
$ws.Cells.Item(168, 12).Value = 5500
$ws.Cells.Item(168, 13).Value = 5500
$ws.Cells.Item(168, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(168, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(168, 16).Value = 275
$ws.Cells.Item(168, 17).Value = 20
$ws.Cells.Item(168, 18).Value = "Hortaliza"
